$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card13")

# Correct the card number in row 2 (was mistakenly "2", should be "13" like the sheet)
$ws.Range("A2").Value = "13"

# Fill every empty data cell in columns D:O (rows 2-13) with the literal "nan"
# placeholder text used throughout this workbook for missing values.
$cols = @("D","E","F","G","H","I","J","K","L","M","N","O")
for ($row = 2; $row -le 13; $row++) {
    foreach ($col in $cols) {
        $cell = $ws.Range("$col$row")
        $cur = $cell.Value2
        if ($cur -eq $null -or $cur -eq "") {
            $cell.Value = "nan"
        }
    }
}
